$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a third training session block (rows 13-18), mirroring the layout of
# the "Session 2" block (rows 7-12). Training for this session was stopped
# early, so epochs 11 and 12 have no loss values (they show as #DIV/0!).
# ---------------------------------------------------------------------------

# Row 13: "Session 3" header
$ws.Range("A13").Value = "Session 3"

# Row 14: column headers
$ws.Range("A14").Value = "Epoch"
$ws.Range("B14").Value = "Training Loss"
$ws.Range("C14").Value = "Validation Loss"
$ws.Range("D14").Value = "Improvement"
$ws.Range("E14").Value = "Improvement"

# Row 15: epoch 9 data (first row of session, no improvement yet)
$ws.Range("A15").Value = 9
$ws.Range("B15").Value = 0.13750000000000001
$ws.Range("C15").Value = 0.19826299999999999
$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "-"

# Row 16: epoch 10 data
$ws.Range("A16").Value = 10
$ws.Range("B16").Value = 0.1052
$ws.Range("C16").Value = 0.12718599999999999
$ws.Range("D16").Formula = "=(B15-B16)/B15"
$ws.Range("E16").Formula = "=(C15-C16)/C15"

# Row 17: epoch 11 - training stopped, loss cells blank
$ws.Range("A17").Value = 11
$ws.Range("D17").Formula = "=(B16-B17)/B16"
$ws.Range("E17").Formula = "=(C16-C17)/C16"

# Row 18: epoch 12 - training stopped, loss cells blank
$ws.Range("A18").Value = 12
$ws.Range("D18").Formula = "=(B17-B18)/B17"
$ws.Range("E18").Formula = "=(C17-C18)/C17"

# Copy the formatting of the Session 2 block onto the new Session 3 block.
$ws.Range("A7:E12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

# Merge the new session header row, same as the other session header rows.
$ws.Range("A13:E13").Merge()

# ---------------------------------------------------------------------------
# Conditional formatting: highlight Validation Loss (column C) relative to
# Training Loss (B3) for each of the three sessions.
# ---------------------------------------------------------------------------

# Session 1 (C3:C6): green when greater than B3, red when less than B3.
$rng1 = $ws.Range("C3:C6")
$cf1a = $rng1.FormatConditions.Add(1, 5, "=`$B`$3")
$cf1a.Font.Color = 13561798
$cf1a.Interior.Color = 10283381
$cf1b = $rng1.FormatConditions.Add(1, 6, "=`$B`$3")
$cf1b.Font.Color = 255
$cf1b.Interior.Color = 12189695

# Session 2 (C9:C12): red when less than B3, green when greater than B3.
$rng2 = $ws.Range("C9:C12")
$cf2a = $rng2.FormatConditions.Add(1, 6, "=`$B`$3")
$cf2a.Font.Color = 255
$cf2a.Interior.Color = 12189695
$cf2b = $rng2.FormatConditions.Add(1, 5, "=`$B`$3")
$cf2b.Font.Color = 13561798
$cf2b.Interior.Color = 10283381

# Session 3 (C15:C18): red when less than B3, green when greater than B3.
$rng3 = $ws.Range("C15:C18")
$cf3a = $rng3.FormatConditions.Add(1, 6, "=`$B`$3")
$cf3a.Font.Color = 255
$cf3a.Interior.Color = 12189695
$cf3b = $rng3.FormatConditions.Add(1, 5, "=`$B`$3")
$cf3b.Font.Color = 13561798
$cf3b.Interior.Color = 10283381

# Keep the visible selection on the newly-added session block, like the
# previous selection highlighted the Session 2 block.
$ws.Range("A14:C16").Select()
